# Applies the "Add files via upload" changes to cronogramaSegundaEntrega.xlsx
# (cronograma / schedule: fill in start / agreed-delivery / delivery dates
#  for several tasks, plus assign two more tasks to team members).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "d-mmm"

# --- Row 4: Escopo do Projeto ---------------------------------------------
# I4 had no date before and still has none, but its format becomes the date
# style used elsewhere in the sheet.
$ws.Range("I4").NumberFormat = $dateFmt

# --- Row 5: Modelo Conceitual ----------------------------------------------
$ws.Range("J5").Value = 43535
$ws.Range("J5").NumberFormat = $dateFmt

# --- Row 6: Modelo Lógico ---------------------------------------------------
$ws.Range("I6").Value = 43535
$ws.Range("I6").NumberFormat = $dateFmt
$ws.Range("K6").Value = 43545
$ws.Range("K6").NumberFormat = $dateFmt

# --- Row 7: Dicionário de dados --------------------------------------------
$ws.Range("I7").Value = 43545
$ws.Range("I7").NumberFormat = $dateFmt
$ws.Range("K7").Value = 43545
$ws.Range("K7").NumberFormat = $dateFmt

# --- Row 8: Wireframe web ---------------------------------------------------
$ws.Range("I8").Value = 43531
$ws.Range("I8").NumberFormat = $dateFmt
$ws.Range("J8").Value = 43535
$ws.Range("J8").NumberFormat = $dateFmt

# --- Row 9: Wireframe mobile -------------------------------------------------
$ws.Range("I9").Value = 43539
$ws.Range("I9").NumberFormat = $dateFmt
$ws.Range("J9").Value = 43541
$ws.Range("J9").NumberFormat = $dateFmt
$ws.Range("K9").Value = 43541
$ws.Range("K9").NumberFormat = $dateFmt

# --- Row 10: Wireframe desktop ----------------------------------------------
$ws.Range("K10").Value = 43542
$ws.Range("K10").NumberFormat = $dateFmt

# --- Row 11: Wireframe CMS ---------------------------------------------------
$ws.Range("I11").Value = 43546
$ws.Range("I11").NumberFormat = $dateFmt
$ws.Range("J11").Value = 43547
$ws.Range("J11").NumberFormat = $dateFmt

# --- Row 12: Frontend site ---------------------------------------------------
$ws.Range("I12").Value = 43543
$ws.Range("I12").NumberFormat = $dateFmt

# --- Row 13: Frontend CMS ----------------------------------------------------
$ws.Range("I13").Value = 43539
$ws.Range("I13").NumberFormat = $dateFmt
# K13 stays empty but switches to the date-formatted style.
$ws.Range("K13").NumberFormat = $dateFmt

# --- Row 14: 2 CRUD CMS -------------------------------------------------------
$ws.Range("H14").Value = "Kaio, Igor"
$ws.Range("I14").Value = 43547
$ws.Range("I14").NumberFormat = $dateFmt
$ws.Range("J14").Value = 43548
$ws.Range("J14").NumberFormat = $dateFmt

# --- Row 15: UML ---------------------------------------------------------------
$ws.Range("H15").Value = "Sarah, Manu, Leonardo"

# Resize column H (Responsáveis) to fit the newly added, longer names.
$ws.Columns.Item(8).AutoFit()

# Leave the active selection on I5, matching where the author was working.
$ws.Range("I5").Select()
